# Update "docs/epexspot_prices.xlsx" with the latest day of data (28-aug)
# across the three sheets: "Prix Spot" (hourly spot prices), "Gaz" (gas
# daily closing price) and "CO2" (CO2 daily closing price).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Prix Spot": append a new date column BX (28-aug) with its 24
# hourly values, copying the header formatting from the previous date
# column (BW) so the new header cell matches the existing style.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Prix Spot")

$ws1.Range("BX1").Value = "28-aug"
$ws1.Range("BW1").Copy() | Out-Null
$ws1.Range("BX1").PasteSpecial(-4122) | Out-Null

$ws1.Range("BX2").Value = 97.15000000000001
$ws1.Range("BX3").Value = 90.14
$ws1.Range("BX4").Value = 86.06999999999999
$ws1.Range("BX5").Value = 73.25
$ws1.Range("BX6").Value = 69.40000000000001
$ws1.Range("BX7").Value = 76.42
$ws1.Range("BX8").Value = 83.38
$ws1.Range("BX9").Value = 90.38
$ws1.Range("BX10").Value = 90.38
$ws1.Range("BX11").Value = 69.40000000000001
$ws1.Range("BX12").Value = 42.57
$ws1.Range("BX13").Value = 20.46
$ws1.Range("BX14").Value = 17.21
$ws1.Range("BX15").Value = 2.09
$ws1.Range("BX16").Value = 0.65
$ws1.Range("BX17").Value = 0.01
$ws1.Range("BX18").Value = 0.65
$ws1.Range("BX19").Value = 21.04
$ws1.Range("BX20").Value = 42.46
$ws1.Range("BX21").Value = 71.12
$ws1.Range("BX22").Value = 90
$ws1.Range("BX23").Value = 97.08
$ws1.Range("BX24").Value = 93.77
$ws1.Range("BX25").Value = 88.09

# ---------------------------------------------------------------------
# Sheet "Gaz": append row 73 with the new day's closing price. The date
# is stored as plain text (like every other date in column A), so the
# cell is briefly formatted as Text to stop Excel from auto-converting
# the "2025-08-26" entry into a date serial number, then the formatting
# is cleared again so the cell ends up unstyled, same as its neighbours.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Gaz")
$ws2.Range("A73").NumberFormat = "@"
$ws2.Range("A73").Value = "2025-08-26"
$ws2.Range("A73").ClearFormats()
$ws2.Range("B73").Value = 32.175

# ---------------------------------------------------------------------
# Sheet "CO2": append row 73 with the new day's closing price (same
# text-date handling as above).
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("CO2")
$ws3.Range("A73").NumberFormat = "@"
$ws3.Range("A73").Value = "2025-08-26"
$ws3.Range("A73").ClearFormats()
$ws3.Range("B73").Value = 72.28
